$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns from generic "Characteristic N" labels to the
# Russian department/type labels used by the app.
$ws.Range("B1").Value = "Отдел"
$ws.Range("C1").Value = "Тип"

# Move the active selection to the last data row (C7) instead of the
# stale C11 reference left over from a larger sample sheet.
$ws.Range("C7").Select()
